# Autogenerated on Wed Apr 01 2015 00:15:40 GMT+0000 (Coordinated Universal Time)
# Refresh MSME summary percentages with more precise (2-decimal) figures.
# The cells hold text (not numbers), so the sheet's General-format auto
# number detection is avoided by forcing a text format before/after the
# write - this keeps the values stored as plain text, matching the source.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Enterprises density (per 1000 people) - row 11
Set-TextValue "B11" "20.92"
Set-TextValue "C11" "0.55"
Set-TextValue "D11" "21.46"

# Employment (% of total) - row 12
Set-TextValue "B12" "70.72"
Set-TextValue "C12" "15.73"
Set-TextValue "D12" "86.45"

# Enterprises (% of total) - row 14
Set-TextValue "B14" "97.33"
Set-TextValue "C14" "2.55"
Set-TextValue "D14" "99.88"
